# Updated cryptos list (Price + Volume(1h) columns) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.338.59"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "2.037.70"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.75"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.656"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.91"
$ws.Range("E8").Value = "  -7.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.22"
$ws.Range("E9").Value = "  +5.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.358"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0737"
$ws.Range("E11").Value = "  -4.29%  "
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.932"
$ws.Range("E13").Value = "  +7.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.37"
$ws.Range("E14").Value = "  -4.39%  "
$ws.Range("D15").Value = "2.334.34"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.32"
$ws.Range("E16").Value = "  -4.83%  "
$ws.Range("D17").Value = "2.047.32"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").Value = "36.287.66"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.81"
$ws.Range("E19").Value = "  -5.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.92"
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("D21").Value = "0.0₃0842"
$ws.Range("E21").Value = "  -4.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.72"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.12"
$ws.Range("E23").Value = "  -4.23%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -3.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.25"
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.07"
$ws.Range("E28").Value = "  -11.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.66"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.120"
$ws.Range("E30").Value = "  -2.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +6.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.91"
$ws.Range("E32").Value = "  -10.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0586"
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.34"
$ws.Range("E34").Value = "  -7.22%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0852"
$ws.Range("E36").Value = "  +4.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.19"
$ws.Range("E38").Value = "  -4.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.90"
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("E40").Value = "  -6.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.89"
$ws.Range("E41").Value = "  -4.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0211"
$ws.Range("E42").Value = "  -4.67%  "
$ws.Range("E43").Value = "  -4.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "92.16"
$ws.Range("E44").Value = "  -3.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0885"
$ws.Range("E45").Value = "  -5.35%  "
$ws.Range("D46").Value = "1.372.25"
$ws.Range("E46").Value = "  +5.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.58"
$ws.Range("E47").Value = "  -6.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.30"
$ws.Range("E48").Value = "  +9.17%  "
$ws.Range("E49").Value = "  +2.10%  "
$ws.Range("D50").Value = "2.223.37"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("E51").Value = "  -4.51%  "
